$d = $word.ActiveDocument

# The existing "_GoBack" bookmark sits on the first paragraph; it moves
# to the end of the new second paragraph, so drop it here first.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Collapse to the very end of the document (end of paragraph 1, right
# before the sectPr) and insert the new, second paragraph as raw OOXML
# so we get the split runs / proofErr spell-check markers / bookmark
# exactly as Word would have produced them.
$end = $d.Content
$end.Collapse(0)

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
       '<pkg:xmlData>' +
       '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:body>' +
       '<w:p>' +
       '<w:pPr><w:jc w:val="center"/></w:pPr>' +
       '<w:proofErr w:type="spellStart"/>' +
       '<w:r><w:t>Commit</w:t></w:r>' +
       '<w:proofErr w:type="spellEnd"/>' +
       '<w:r><w:t xml:space="preserve"> inicial do meu livro</w:t></w:r>' +
       '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
       '<w:bookmarkEnd w:id="0"/>' +
       '</w:p>' +
       '</w:body>' +
       '</pkg:xmlData>' +
       '</pkg:part>' +
       '</pkg:package>'

$end.InsertXML($xml)
